$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B (target OOXML width 15.42578125, was 14.42578125).
# The COM layer quantizes ColumnWidth to steps of 1/6 when it writes the
# OOXML "width" attribute, so 14.666666666666666 is the closest input that
# reliably rounds to the nearest reachable width (15.5).
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# Update existing values in rows 1-4
$ws.Range("A1").Value = -0.0044275530382215976
$ws.Range("B1").Value = -0.0041450133963361873

$ws.Range("A2").Value = -0.018291144398067911
$ws.Range("B2").Value = -0.021493930888234463

$ws.Range("A3").Value = -0.037334577392695828
$ws.Range("B3").Value = -0.036950788652470062

$ws.Range("A4").Value = -0.0049109377680854101
$ws.Range("B4").Value = -0.0049106643264613253

# Add new row 5
$ws.Range("A5").Value = -0.070935905377017966
$ws.Range("B5").Value = -0.070934448574766482
